# Applies the "Add files via upload" commit:
#   - Sprint 2 sheet gets task rows for US14 (Multiple births <= 5) and
#     US16 (Male last names), each with owner/status/estimate/actual.
#   - Backlog sheet gets Owner/Status for those two backlog items.
#   - Sprint1 (Sprint 1 review) sheet gets "Keep doing:" / "Avoid:"
#     retrospective comments from ALW.
#   - Various view-state tweaks (selected cell, active tab, column widths).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Backlog: mark US14 / US16 backlog rows with an owner + status
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

$backlog.Range("D10").Value = "ALW"
$backlog.Range("E10").Value = "Not started"
$backlog.Range("D11").Value = "ALW"
$backlog.Range("E11").Value = "Not started"

# Column C / E grow to fit the new longer text (closest value this
# engine's column-width grid can represent).
$backlog.Columns.Item(3).ColumnWidth = 25.714285714285715
$backlog.Columns.Item(5).ColumnWidth = 9.714285714285714

# ---------------------------------------------------------------------
# Sprint1: Sprint 1 retrospective notes ("Keep doing:" / "Avoid:")
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

$sprint1.Range("H8").Value = 60

function Set-Note($ws, $cellRef, $text, $height) {
    $c = $ws.Range($cellRef)
    $c.Value = $text
    $c.NumberFormat = "@"
    $c.WrapText = $true
    $c.Font.Name = "Verdana"
    $c.Font.Size = 10
    $c.Font.Bold = $false
    $row = $cellRef -replace '[A-Za-z]', ''
    $ws.Rows.Item([int]$row).RowHeight = $height
}

# "Keep doing:" block (row 37) grows with two new bullet rows
Set-Note $sprint1 "B38" "From ALW: GitHub is working great" 25.5
Set-Note $sprint1 "B39" "From ALW: Communication between the team members, appears to be working as far as I'm concerned" 63.75

# "Avoid:" block (row 41) grows with three new bullet rows
Set-Note $sprint1 "B42" "From ALW: I was too liberal with my estimated time. The time spent took longer than I initially documented. I need to be more conserative with my time spent" 89.25
Set-Note $sprint1 "B43" "From ALW: Adding unit testing into the tasks list. This was additional work that wasn't created initially" 51
Set-Note $sprint1 "B44" "From ALW: Programming wise, needed to break my US into separate functions, so unit testing could be easier" 63.75

# ---------------------------------------------------------------------
# Sprint2: populate the Sprint-2 task breakdown for US14 and US16
# ---------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint2")

$sprint2.Columns.Item(2).ColumnWidth = 17.285714285714285

# US14 - Multiple births <= 5
$sprint2.Range("A2").Value = "US14"
$sprint2.Range("A2").Font.Bold = $true
$sprint2.Range("B2").Formula = "=Backlog!`$C`$10"
$sprint2.Range("C2").Value = "ALW"
$sprint2.Range("D2").Value = "Not Started"
$sprint2.Range("E2").Value = 60
$sprint2.Range("F2").Value = 60

$sprint2.Range("A3").Value = "T14.01"
$sprint2.Range("B3").Value = "store children with the same birth date in a collection"
$sprint2.Range("B3").WrapText = $true
$sprint2.Rows.Item(3).RowHeight = 38.25

$sprint2.Range("A4").Value = "T14.02"
$sprint2.Range("B4").Value = "Collect size, throw error if > 5"
$sprint2.Range("B4").WrapText = $true
$sprint2.Rows.Item(4).RowHeight = 25.5

$sprint2.Range("A5").Value = "T14.03"
$sprint2.Range("B5").Value = "Create Unit Test to test"
$sprint2.Range("B5").WrapText = $true
$sprint2.Rows.Item(5).RowHeight = 25.5

$sprint2.Range("A6").Value = ""
$sprint2.Range("B6").Value = ""

# US16 - Male last names
$sprint2.Range("A7").Value = "US16"
$sprint2.Range("A7").Font.Bold = $true
$sprint2.Range("B7").Formula = "=Backlog!`$C`$11"
$sprint2.Range("C7").Value = "ALW"
$sprint2.Range("D7").Value = "Not Started"
$sprint2.Range("E7").Value = 60
$sprint2.Range("F7").Value = 60

$sprint2.Range("A8").Value = "T16.01"
$sprint2.Range("B8").Value = "Save the fathers last name"
$sprint2.Range("B8").WrapText = $true
$sprint2.Rows.Item(8).RowHeight = 25.5

$sprint2.Range("A9").Value = "T16.02"
$sprint2.Range("B9").Value = "Compare all males in the family have the same last name"
$sprint2.Range("B9").WrapText = $true
$sprint2.Rows.Item(9).RowHeight = 38.25

$sprint2.Range("A10").Value = "T16.03"
$sprint2.Range("B10").Value = "Create Unit Test to test"
$sprint2.Range("B10").WrapText = $true
$sprint2.Rows.Item(10).RowHeight = 25.5

$sprint2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# View state: selections per sheet, then activate Sprint2 last so it
# ends up the active tab (matches the workbook's activeTab bump).
# ---------------------------------------------------------------------
$stories = $wb.Worksheets.Item("Stories")
$stories.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 13
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$stories.Range("C5").Select()

$backlog.Range("C11").Select()

$sprint1.Range("B40").Select()
try {
    $excel.ActiveWindow.ScrollRow = 35
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$sprint2.Range("E8").Select()
